$d = $word.ActiveDocument

# 1. Title: "Maksim Kanev's Technical Test" -> "Maksim Kanev's Technical Challenge"
$d.Content.Find.Execute("Technical Test", $false, $true, $false, $false, $false, $true, 1, $false, "Technical Challenge", 2)

# 2. Merge the split runs in the three "Objective" paragraphs (text itself is unchanged,
#    but the run boundary that fell in the middle of a sentence is removed).
$d.Content.Find.Execute("should be done on the last day of the month.", $false, $true, $false, $false, $false, $true, 1, $false, "should be done on the last day of the month.", 2)
$d.Content.Find.Execute("should be arranged for the following Monday.", $false, $true, $false, $false, $false, $true, 1, $false, "should be arranged for the following Monday.", 2)
$d.Content.Find.Execute("testing should be set for the previous Thursday.", $false, $true, $false, $false, $false, $true, 1, $false, "testing should be set for the previous Thursday.", 2)

# 3. Fix the "localohst" typo -> "localhost"
$d.Content.Find.Execute("localohst", $false, $true, $false, $false, $false, $true, 1, $false, "localhost", 2)

# 4. Move the "_GoBack" bookmark from the end of the last paragraph to inside
#    "localhost" (right after "localho", before "st"), matching where the author's
#    cursor ended up after fixing the typo.
$full = $d.Content.Text
$idx = $full.IndexOf("localhost")
if ($idx -ge 0) {
    $pos = $idx + 7
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
